{"js": "// Replace the date line and every \"NN\u00d7NN=\" problem cell with its new value.\n// Each old string is unique in the document, so an exact, case-sensitive\n// whole-text search safely targets exactly one run each.\nconst replacements = [\n  [\"2025-06-15 Sunday\", \"2025-06-16 Monday\"],\n  [\"26\u00d778=\", \"47\u00d790=\"],\n  [\"59\u00d792=\", \"12\u00d781=\"],\n  [\"73\u00d743=\", \"24\u00d762=\"],\n  [\"48\u00d722=\", \"73\u00d778=\"],\n  [\"18\u00d782=\", \"89\u00d772=\"],\n  [\"28\u00d782=\", \"77\u00d790=\"],\n  [\"47\u00d751=\", \"80\u00d759=\"],\n  [\"17\u00d778=\", \"85\u00d743=\"],\n  [\"15\u00d711=\", \"15\u00d746=\"],\n  [\"48\u00d768=\", \"14\u00d752=\"],\n  [\"49\u00d773=\", \"63\u00d728=\"],\n  [\"51\u00d735=\", \"18\u00d799=\"],\n  [\"71\u00d790=\", \"92\u00d711=\"],\n  [\"16\u00d740=\", \"40\u00d756=\"],\n  [\"34\u00d748=\", \"68\u00d735=\"],\n  [\"52\u00d757=\", \"77\u00d793=\"],\n  [\"74\u00d771=\", \"37\u00d785=\"],\n  [\"91\u00d751=\", \"29\u00d798=\"],\n  [\"49\u00d769=\", \"78\u00d774=\"],\n  [\"18\u00d798=\", \"15\u00d729=\"],\n  [\"80\u00d712=\", \"26\u00d769=\"],\n  [\"84\u00d789=\", \"34\u00d774=\"],\n  [\"69\u00d795=\", \"26\u00d769=\"],\n  [\"79\u00d711=\", \"80\u00d760=\"],\n  [\"23\u00d766=\", \"76\u00d737=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"NN\u00d7NN=\" problem cell with its new value.\n# Each old string is unique in the document, so Find/Execute with\n# MatchWholeWord/MatchWildcards off and Replace:=wdReplaceAll (2) safely\n# rewrites exactly one occurrence per pair without touching anything else.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-15 Sunday\", \"2025-06-16 Monday\"),\n    @(\"26\u00d778=\", \"47\u00d790=\"),\n    @(\"59\u00d792=\", \"12\u00d781=\"),\n    @(\"73\u00d743=\", \"24\u00d762=\"),\n    @(\"48\u00d722=\", \"73\u00d778=\"),\n    @(\"18\u00d782=\", \"89\u00d772=\"),\n    @(\"28\u00d782=\", \"77\u00d790=\"),\n    @(\"47\u00d751=\", \"80\u00d759=\"),\n    @(\"17\u00d778=\", \"85\u00d743=\"),\n    @(\"15\u00d711=\", \"15\u00d746=\"),\n    @(\"48\u00d768=\", \"14\u00d752=\"),\n    @(\"49\u00d773=\", \"63\u00d728=\"),\n    @(\"51\u00d735=\", \"18\u00d799=\"),\n    @(\"71\u00d790=\", \"92\u00d711=\"),\n    @(\"16\u00d740=\", \"40\u00d756=\"),\n    @(\"34\u00d748=\", \"68\u00d735=\"),\n    @(\"52\u00d757=\", \"77\u00d793=\"),\n    @(\"74\u00d771=\", \"37\u00d785=\"),\n    @(\"91\u00d751=\", \"29\u00d798=\"),\n    @(\"49\u00d769=\", \"78\u00d774=\"),\n    @(\"18\u00d798=\", \"15\u00d729=\"),\n    @(\"80\u00d712=\", \"26\u00d769=\"),\n    @(\"84\u00d789=\", \"34\u00d774=\"),\n    @(\"69\u00d795=\", \"26\u00d769=\"),\n    @(\"79\u00d711=\", \"80\u00d760=\"),\n    @(\"23\u00d766=\", \"76\u00d737=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
